$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, regardless of whether it looks like
# a number (e.g. "215.67") or a percent-ish string (e.g. "  +3.44%  ").
# Forcing NumberFormat to Text before the write keeps Excel from
# reinterpreting the literal as a number, then the style is reset back to
# Normal/default so we don't leave a stray style applied to the cell.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Column D (price) and E (volume 1h) updates, keyed by row number.
# Only rows/cells that actually change are listed here.

$updates = @{
    2  = @{ D = "27.137.85"; E = "  +3.44%  " }
    3  = @{ D = "1.661.27";  E = "  +4.04%  " }
    4  = @{                  E = "  -0.01%  " }
    5  = @{ D = "215.67";    E = "  +1.75%  " }
    6  = @{                  E = "  +1.08%  " }
    7  = @{                  E = "  -0.05%  " }
    8  = @{                  E = "  +2.44%  " }
    9  = @{                  E = "  +1.58%  " }
    10 = @{ D = "19.66";     E = "  +3.79%  " }
    11 = @{ D = "0.0864";    E = "  +1.20%  " }
    12 = @{                  E = "  +3.92%  " }
    13 = @{ D = "1.686.70";  E = "  +5.10%  " }
    14 = @{                  E = "  +1.95%  " }
    15 = @{ D = "0.521";     E = "  +3.10%  " }
    16 = @{                  E = "  +2.13%  " }
    17 = @{ D = "240.84";    E = "  +5.73%  " }
    18 = @{ D = "27.123.21"; E = "  +3.42%  " }
    19 = @{ D = "7.86";      E = "  +3.97%  " }
    21 = @{                  E = "  -0.08%  " }
    22 = @{                  E = "  +5.27%  " }
    23 = @{                  E = "  +3.74%  " }
    24 = @{ D = "9.33";      E = "  +4.87%  " }
    25 = @{ D = "145.76";    E = "  +0.16%  " }
    26 = @{                  E = "  +0.00%  " }
    27 = @{ D = "7.17" }
    28 = @{                  E = "  +1.21%  " }
    29 = @{                  E = "  +3.33%  " }
    30 = @{                  E = "  +1.19%  " }
    31 = @{                  E = "  +1.41%  " }
    32 = @{ D = "1.534.74" }
    33 = @{                  E = "  +3.03%  " }
    34 = @{                  E = "  +3.85%  " }
    35 = @{                  E = "  +8.60%  " }
    36 = @{                  E = "  -0.02%  " }
    37 = @{ D = "0.576";     E = "  +2.11%  " }
    38 = @{                  E = "  +9.37%  " }
    39 = @{                  E = "  +3.14%  " }
    40 = @{                  E = "  +3.89%  " }
    41 = @{                  E = "  -0.06%  " }
    42 = @{ D = "2.27";      E = "  +4.35%  " }
    43 = @{ D = "66.17";     E = "  +9.53%  " }
    44 = @{ D = "1.800.34";  E = "  +3.78%  " }
    45 = @{ D = "0.773";     E = "  +2.11%  " }
    46 = @{ D = "0.918";     E = "  -1.10%  " }
    47 = @{ D = "90.47";     E = "  +3.34%  " }
    48 = @{                  E = "  +4.18%  " }
    49 = @{                  E = "  -0.34%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        Set-TextValue $ws.Range("D$row") $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        Set-TextValue $ws.Range("E$row") $vals["E"]
    }
}

# Rows 50 and 51: Algorand and Cronos swap places (Cronos now ranks
# above Algorand), each with refreshed price/volume figures.
Set-TextValue $ws.Range("B50") "Cronos"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.0504"
Set-TextValue $ws.Range("E50") "  +0.86%  "

Set-TextValue $ws.Range("B51") "Algorand"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D51") "0.0980"
Set-TextValue $ws.Range("E51") "  +3.56%  "
